$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.667.10'
$ws.Range('E2').Value = '  +1.26%  '
$ws.Range('D3').Value = '2.287.44'
$ws.Range('E3').Value = '  +1.56%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '504.36'
$ws.Range('E5').Value = '  +2.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.29'
$ws.Range('E6').Value = '  +2.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.995'
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.529'
$ws.Range('E8').Value = '  +0.88%  '
$ws.Range('D9').Value = '2.300.87'
$ws.Range('E9').Value = '  +1.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0961'
$ws.Range('E10').Value = '  +1.26%  '
$ws.Range('E11').Value = '  +1.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.338'
$ws.Range('E12').Value = '  +4.36%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.93'
$ws.Range('E13').Value = '  +6.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.17'
$ws.Range('E14').Value = '  +6.75%  '
$ws.Range('D15').Value = '2.692.57'
$ws.Range('E15').Value = '  +1.57%  '
$ws.Range('D16').Value = '54.636.36'
$ws.Range('E16').Value = '  +1.36%  '
$ws.Range('E17').Value = '  +1.88%  '
$ws.Range('D18').Value = '2.279.24'
$ws.Range('E18').Value = '  +1.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.35'
$ws.Range('E19').Value = '  +3.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.17'
$ws.Range('E20').Value = '  +2.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '306.90'
$ws.Range('E21').Value = '  +2.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.35'
$ws.Range('E22').Value = '  -0.93%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.69'
$ws.Range('E24').Value = '  -1.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.993'
$ws.Range('E25').Value = '  -2.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.151'
$ws.Range('E26').Value = '  +2.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.41'
$ws.Range('E27').Value = '  +5.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '172.54'
$ws.Range('E28').Value = '  +3.78%  '
$ws.Range('D29').Value = '0.0₃0709'
$ws.Range('E29').Value = '  +4.82%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.63'
$ws.Range('E30').Value = '  +2.08%  '
$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.06'
$ws.Range('E31').Value = '  +3.77%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.11'
$ws.Range('E32').Value = '  +3.69%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.00'
$ws.Range('E34').Value = '  +2.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.995'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.941'
$ws.Range('E36').Value = '  +6.36%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.21'
$ws.Range('E37').Value = '  +2.79%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.77'
$ws.Range('E38').Value = '  +3.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.26'
$ws.Range('E39').Value = '  +1.65%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.377'
$ws.Range('E40').Value = '  +1.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.42'
$ws.Range('E41').Value = '  +2.35%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.42'
$ws.Range('E42').Value = '  +2.51%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.95'
$ws.Range('E43').Value = '  +0.92%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '125.67'
$ws.Range('E44').Value = '  +0.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0494'
$ws.Range('E45').Value = '  +2.84%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '248.24'
$ws.Range('E46').Value = '  +5.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0899'
$ws.Range('E47').Value = '  +1.54%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.550'
$ws.Range('E48').Value = '  +2.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.376'
$ws.Range('E49').Value = '  +1.92%  '
$ws.Range('E50').Value = '  +3.36%  '
